$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.884.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5107"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2562"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06335"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07780"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.644.06"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.253"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.856.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5499"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.71"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7632"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.905.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.410"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.848"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.018"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.893"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.08"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1251"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.51%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.752"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.240"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04891"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.230"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.541"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.374"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8966"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.114.05"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.573"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7955"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.764.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.53%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.60"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05129"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.548"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.69%  "
